$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("04 Sep")

# Update data values
$ws.Range("E11").Value = 252
$ws.Range("E15").Value = 41959

# Update the selected cell / active cell on the sheet
$ws.Activate()
$ws.Range("E16").Select()
